# Add a new row (33) of data to Sheet1: date 46002 (12/11/2025) -> count 10
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A33").Value = 46002
$ws.Range("A33").NumberFormat = "m/d/yy"
$ws.Range("B33").Value = 10

# Match the saved selection: A33:B33 with A33 as the active cell
$ws.Range("A33:B33").Select() | Out-Null
